$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number for every data row (2-233).
# All of these cells currently contain 45203 and need to be bumped to 45204.
$ws.Range("C2:C233").Value = 45204
